# trying url for imports
#
# This script reproduces the commit that adds several new pelvis-related
# terms to the "axis" and "structures" sheets, adds a new mapped column
# (I) entry on the "trait" sheet, removes the now-redundant
# "shaft of ilium circumference" row from "axis", and leaves the "axis"
# sheet as the active tab/selection instead of "AB".

$wb = $excel.ActiveWorkbook

$wsTrait = $wb.Worksheets.Item("trait")
$wsAxis = $wb.Worksheets.Item("axis")
$wsAB = $wb.Worksheets.Item("AB")
$wsStructures = $wb.Worksheets.Item("structures")

# ---------------------------------------------------------------------
# axis sheet: drop the old row 7 ("shaft of ilium circumference" /
# no axis / no structure) -- the rows below it shift up one place.
# ---------------------------------------------------------------------
$wsAxis.Rows.Item(7).Delete()

# ---------------------------------------------------------------------
# Populate the new cells in the order that matches the commit's shared
# string table (new unique strings are appended in first-use order).
# ---------------------------------------------------------------------
$wsAxis.Range("E2").Value = "This includes the three bones of the pelvis: the pubis, ischium, and ilium"

$wsAxis.Range("C4").Value = "pubic symphsis"
$wsAxis.Range("B4").Value = "proximal-distal"

$wsStructures.Range("B3").Value = "ilium shaft"

$wsAxis.Range("C7").Value = "obtruartor foramen"
$wsAxis.Range("B7").Value = "proximal-distal"

$wsStructures.Range("B4").Value = "tubera coxarum"

$wsStructures.Range("B5").Value = "tubera ischiadica"

$wsAxis.Range("B6").Value = "medial-lateral"

$wsTrait.Range("I9").Value = "circumfernce and ('inheres in' some 'shaft of ilium')"

# ---------------------------------------------------------------------
# Remaining reused-string cells on axis (all "medial-lateral" /
# "acetabular part of hip bone" / "pubic symphsis" duplicates).
# ---------------------------------------------------------------------
$wsAxis.Range("B8").Value = "medial-lateral"
$wsAxis.Range("B9").Value = "medial-lateral"
$wsAxis.Range("C9").Value = "acetabular part of hip bone"
$wsAxis.Range("B10").Value = "medial-lateral"
$wsAxis.Range("B11").Value = "medial-lateral"

# structures sheet: second new row (pubic symphsis), matching B2
$wsStructures.Range("B2").Value = "pubic symphsis"

# ---------------------------------------------------------------------
# Page setup (both trait and axis gained an explicit portrait pageSetup
# element in the saved file).
# ---------------------------------------------------------------------
$wsTrait.PageSetup.Orientation = 1
$wsAxis.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selections / active tab. AB was active before; axis becomes active
# (with B12 selected) and AB reverts to its own B3 selection without the
# tab being the active one. trait keeps I9 selected, structures keeps B5
# selected, but neither is the active tab -- axis (selected last) is.
# ---------------------------------------------------------------------
[void]$wsAB.Range("B3").Select()
[void]$wsTrait.Range("I9").Select()
[void]$wsStructures.Range("B5").Select()
[void]$wsAxis.Range("B12").Select()
